$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 updates
$ws.Range("B3").Value = "Sandeep"
$ws.Range("C3").Value = "--"

# Row 4 updates
$ws.Range("B4").Value = "Ajay"
$ws.Range("C4").Value = "1Member 1M4Member 4MFMEmber five+ 1M1Member 1 C by OwnerM4Member 4 C by admin 1MFMEmber five created by Admin oneTTtest test"

# Row 5 updates
$ws.Range("B5").Value = "Team 1 created by Admin 1"
$ws.Range("C5").Value = "4Member 4MFMEmber five"

# New row 6 - seed A6/D6 from an existing blank cell so they materialize as
# real (empty) cells rather than being left absent, then fill B6/C6.
$ws.Range("A2").Copy($ws.Range("A6"))
$ws.Range("D2").Copy($ws.Range("D6"))
$ws.Range("B6").Value = "Team 3 edited by owner 2"
$ws.Range("C6").Value = "3Admin 3M4Member 4MFMEmber five"
